# Insert a new data row at row 492 (pushing existing rows 492.. down by one)
# and populate it with the new record's values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 492; this shifts row 492..590 down to 493..591
$ws.Rows.Item(492).Insert()

# Match the date number format used by the rest of column D for the new row.
$ws.Cells.Item(492, 4).NumberFormat = $ws.Cells.Item(491, 4).NumberFormat

# Fill in the new row's values
$ws.Cells.Item(492, 1).Value = 3
$ws.Cells.Item(492, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(492, 3).Value = "Coquimbo"
$ws.Cells.Item(492, 4).Value = 44889
$ws.Cells.Item(492, 5).Value = 5
$ws.Cells.Item(492, 6).Value = 100112003
$ws.Cells.Item(492, 7).Value = "Ajo"
$ws.Cells.Item(492, 8).Value = "Chino"
$ws.Cells.Item(492, 9).Value = "Primera"
$ws.Cells.Item(492, 10).Value = 125
$ws.Cells.Item(492, 11).Value = 13500
$ws.Cells.Item(492, 12).Value = 14000
$ws.Cells.Item(492, 13).Value = 13740
$ws.Cells.Item(492, 14).Value = '$/caja 10 kilos'
$ws.Cells.Item(492, 15).Value = "China"
$ws.Cells.Item(492, 16).Value = 1374
$ws.Cells.Item(492, 17).Value = 10
$ws.Cells.Item(492, 18).Value = "Hortaliza"
